$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New conversation rows scraped from the MySQL db, appended after the
# existing data (rows 2-7) for the "Noah" / 13052054965 thread.
$rows = @(
    @("2025-09-19 20:35:14", "Noah", 8450689526, "13052054965", "Hi"),
    @("2025-09-19 21:03:51", "Noah", 8450689526, "13052054965", "Hi"),
    @("2025-09-19 21:10:32", "Noah", 8450689526, "13052054965", "Hi"),
    @("2025-09-19 21:11:51", "Noah", 8450689526, "13052054965", "Hi"),
    @("2025-09-19 21:12:13", "Noah", 8450689526, "13052054965", "Yo")
)

$r = 8
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    # Phone numbers in column D are kept as text (matches the existing
    # "13052054965" text cells in rows 6-7). A leading apostrophe forces
    # Excel to store the numeric-looking string as text instead of
    # auto-converting it to a number; resetting the style afterwards
    # drops the quote-prefix formatting Excel would otherwise tag the
    # cell with, keeping it on the sheet's default (unstyled) look.
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    $ws.Cells.Item($r, 4).Style = "Normal"

    $ws.Cells.Item($r, 5).Value = $row[4]

    # Media / Channel columns are blank text cells (same empty-string
    # pattern used throughout the existing rows), not simply empty/null
    # cells, so use the same apostrophe trick to keep them as text.
    $ws.Cells.Item($r, 6).Value = "'"
    $ws.Cells.Item($r, 6).Style = "Normal"
    $ws.Cells.Item($r, 7).Value = "'"
    $ws.Cells.Item($r, 7).Style = "Normal"

    $r = $r + 1
}
